$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite the model-name labels in column A (rows 2-23) to the new,
# more interpretable names, and make the AIC tables easier to read.
$ws.Range("A2").Value = "site"
$ws.Range("A3").Value = "rugosity + site"
$ws.Range("A4").Value = "year + site"
$ws.Range("A5").Value = "coralcover + site"
$ws.Range("A6").Value = "spongecover + site"
$ws.Range("A7").Value = "rugosity + year + site"
$ws.Range("A8").Value = "spongecover + year + site"
$ws.Range("A9").Value = "coralcover + year + site"
$ws.Range("A10").Value = "year + site + year*site"
$ws.Range("A11").Value = "rugosity + site + site*rugosity"
$ws.Range("A12").Value = "spongecover + site + site*spongecover"
$ws.Range("A13").Value = "coralcover + site + site*coralcover"
$ws.Range("A14").Value = "rugosity + year + year*rugosity"
$ws.Range("A15").Value = "rugosity + year"
$ws.Range("A16").Value = "rugosity"
$ws.Range("A17").Value = "coralcover + year"
$ws.Range("A18").Value = "coralcover + year + year*coralcover"
$ws.Range("A19").Value = "coralcover"
$ws.Range("A20").Value = "spongecover"
$ws.Range("A21").Value = "spongecover + year + year*spongecover"
$ws.Range("A22").Value = "spongecover + year"
$ws.Range("A23").Value = "year"

# Column A now needs to be wider to fit the longer, more descriptive names
# (matches the ~37-character-wide "spongecover + site + site*spongecover"
# style labels now stored there).
$ws.Columns.Item(1).ColumnWidth = 36.333333

# Reset the saved cell selection / active-cell marker in the sheet view.
$ws.Range("A1").Select()
